$wb = $excel.ActiveWorkbook

# --- Sample_Custom_Moves sheet: fill in row 8 (Big Ears) ---
$wsMoves = $wb.Worksheets.Item("Sample_Custom_Moves")
$wsMoves.Range("B8").Value = "+1 to Sense roll when listening to others' conversations"
$wsMoves.Range("A8").Value = "Big Ears"
# Column B text begins with "+" so Excel stores it with a quote-prefix
# (text) style, matching the formatting already used by the other rows
# in this column (e.g. B7). Copy that cell's format onto B8/B9.
$wsMoves.Range("B7").Copy() | Out-Null
$wsMoves.Range("B8").PasteSpecial(-4122) | Out-Null

# --- Bestiary sheet: fill in row 21 (Walking Shark) ---
$wsBestiary = $wb.Worksheets.Item("Bestiary")
$wsBestiary.Range("A21").Value = "Walking Shark"
$wsBestiary.Range("B21").Value = "A shark on two legs"
$wsBestiary.Range("C21").Value = 9
$wsBestiary.Range("D21").Value = 0
$wsBestiary.Range("E21").Value = "Lightning"
$wsBestiary.Range("F21").Value = "Bite: 1d10 damage.  Once bitten, 3 damage every turn until the enemy escapes."

# --- Sample_Custom_Moves sheet: fill in row 9 (Detective) ---
$wsMoves.Range("A9").Value = "Detective"
$wsMoves.Range("B9").Value = "+1 to Search when investigating a crime"
$wsMoves.Range("B7").Copy() | Out-Null
$wsMoves.Range("B9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Update selections to match final author state ---
[void]$wsBestiary.Range("B39").Select()
[void]$wsMoves.Select()
[void]$wsMoves.Range("B10").Select()
